# 20160411 - 006 / running_logs / logs.xlsx
# "Update results 005 - 011"
#
# The five runs timestamped 20160414_* (previously rows 10-14, using the
# "n_iterator: 100" model and a "0 filters: " w/ trailing-space template
# filter label) move up to rows 2-6, and the template-filter label loses its
# trailing space ("0 filters:"). They are followed by five brand new runs
# timestamped 20160415_* (rows 7-11) that use a re-ordered preprocessing
# pipeline and the "n_iterator: 300" model. The eight oldest runs
# (20160412_*) are dropped entirely, shrinking the sheet from A1:J14 to
# A1:J11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- text blocks shared by several rows -------------------------------------------------
$Features        = "12 features: length, #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, %max_digit_skip_0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit"
$Model           = "Neuron Network"

$PreprocessA     = "convert to lower, convert unicode to ascii, trim `"space`" and `",`", space after punctuation, remove multiple spaces, remove break line"
$ModelDetails100 = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 100"
$FilterA         = "0 filters:"

$PreprocessB     = "trim `"space`" and `",`", remove break line, remove multiple spaces, convert unicode to ascii, space after punctuation, convert to lower"
$ModelDetails300 = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300"
$FilterB         = "0 filters: "

# --- the 10 final data rows (A:J), row 1 is the untouched header -----------------------
$rows = @(
    @{ Row=2;  Time="20160414_100653"; Running=210.168;  Preprocess=$PreprocessA; ModelDetails=$ModelDetails100; Filter=$FilterA; Test=0.999333333333333; Val=0.937293729372937; Template=0.0759493670886076 }
    @{ Row=3;  Time="20160414_101023"; Running=234.987;  Preprocess=$PreprocessA; ModelDetails=$ModelDetails100; Filter=$FilterA; Test=1;                  Val=0.933993399339934; Template=0.0769230769230769 }
    @{ Row=4;  Time="20160414_101418"; Running=212.297;  Preprocess=$PreprocessA; ModelDetails=$ModelDetails100; Filter=$FilterA; Test=1;                  Val=0.933993399339934; Template=0.0769230769230769 }
    @{ Row=5;  Time="20160414_101751"; Running=190.864;  Preprocess=$PreprocessA; ModelDetails=$ModelDetails100; Filter=$FilterA; Test=1;                  Val=0.937293729372937; Template=0.0886075949367089 }
    @{ Row=6;  Time="20160414_102102"; Running=208.347;  Preprocess=$PreprocessA; ModelDetails=$ModelDetails100; Filter=$FilterA; Test=1;                  Val=0.933993399339934; Template=0.0897435897435897 }
    @{ Row=7;  Time="20160415_171134"; Running=198.666;  Preprocess=$PreprocessB; ModelDetails=$ModelDetails300; Filter=$FilterB; Test=1;                  Val=0.933993399339934; Template=0.115384615384615  }
    @{ Row=8;  Time="20160415_171453"; Running=252.985;  Preprocess=$PreprocessB; ModelDetails=$ModelDetails300; Filter=$FilterB; Test=1;                  Val=0.933993399339934; Template=0.115384615384615  }
    @{ Row=9;  Time="20160415_171906"; Running=288.741;  Preprocess=$PreprocessB; ModelDetails=$ModelDetails300; Filter=$FilterB; Test=0.998;               Val=0.943894389438944; Template=0.0740740740740741 }
    @{ Row=10; Time="20160415_172354"; Running=314.434;  Preprocess=$PreprocessB; ModelDetails=$ModelDetails300; Filter=$FilterB; Test=1;                  Val=0.933993399339934; Template=0.102564102564103  }
    @{ Row=11; Time="20160415_172909"; Running=346.333;  Preprocess=$PreprocessB; ModelDetails=$ModelDetails300; Filter=$FilterB; Test=1;                  Val=0.933993399339934; Template=0.115384615384615  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Time
    $ws.Range("B$n").Value = $r.Running
    $ws.Range("C$n").Value = $r.Preprocess
    $ws.Range("D$n").Value = $Features
    $ws.Range("E$n").Value = $Model
    $ws.Range("F$n").Value = $r.ModelDetails
    $ws.Range("G$n").Value = $r.Test
    $ws.Range("H$n").Value = $r.Val
    $ws.Range("I$n").Value = $r.Filter
    $ws.Range("J$n").Value = $r.Template
}

# The table used to run through row 14; the three trailing rows are gone now.
$ws.Range("A12:J14").ClearContents()
